# Applies the 11/12/2017 MAMATHA CHICK IN edit:
#  1. Merge the split "SUN Dec 03" / " 09:34:26 PST 2017" runs into one run.
#  2. Append a new purchase-entry block (SUN Dec 10 ...) after the last
#     recorded entry ("Amount Received mode ... - CASH").

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: merge the two runs that make up the "SUN Dec 03 ..." date
# line into a single run/text node.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "SUN Dec 03 09:34:26 PST 2017", $false, $false, $false, $false, $false,
    $true, 1, $false, "SUN Dec 03 09:34:26 PST 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: add the new "SUN Dec 10" purchase-entry block right after the
# last non-empty paragraph in the document (the previous entry's closing
# "Amount Received mode ... - CASH" line), before the trailing blank
# paragraphs.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -gt 0) {
        $targetIndex = $i
        break
    }
}

$lastPara = $d.Paragraphs.Item($targetIndex)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$font = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'
$fontRed = '<w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FF0000"/>'

$xml = "<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr></w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr>" +
"<w:r><w:rPr>$font</w:rPr><w:t>SUN Dec 10</w:t></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:t xml:space=""preserve""> 10:03:22 PST 2017</w:t></w:r>" +
"</w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr>" +
"<w:r><w:rPr>$font</w:rPr><w:t>Person Name</w:t></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/><w:t>- PUTTARAJU</w:t></w:r>" +
"</w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr>" +
"<w:r><w:rPr>$font</w:rPr><w:t>---------------------------------------------------------------</w:t></w:r>" +
"</w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr>" +
"<w:r><w:rPr>$font</w:rPr><w:t>Item Name</w:t></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/><w:t>- POTATO</w:t></w:r>" +
"</w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$fontRed</w:rPr></w:pPr>" +
"<w:r><w:rPr>$fontRed</w:rPr><w:t>Amount Received</w:t></w:r>" +
"<w:r><w:rPr>$fontRed</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$fontRed</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$fontRed</w:rPr><w:tab/><w:t>- 612</w:t></w:r>" +
"</w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr>" +
"<w:r><w:rPr>$font</w:rPr><w:t>Amount Received mode</w:t></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/></w:r>" +
"<w:r><w:rPr>$font</w:rPr><w:tab/><w:t>- CASH AND CLEARD</w:t></w:r>" +
"</w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr></w:p>" +
"<w:p $w><w:pPr><w:pStyle w:val=""PlainText""/><w:rPr>$font</w:rPr></w:pPr></w:p>"

$insertPoint.InsertXML($xml)

Write-Host "Edit applied."
